$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change row 4: "Zusatzbeitrag Krankenversicherung in Prozent" -> AG-Anteil, value stays percent-like but changes to 0.99
$ws.Range("A4").Value = "Zusatzbeitrag Krankenversicherung AG-Anteil in Prozent"
$ws.Range("B4").Value = 0.99

# Insert a new row 5 for AN-Anteil (shifts everything below down by one)
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "Zusatzbeitrag Krankenversicherung AN-Anteil in Prozent"
$ws.Range("B5").Value = 0.99
$ws.Range("B5").NumberFormat = "0.00"

# Update Umlage U1 value (now row 6)
$ws.Range("B6").Value = 2.3

# Row 9 is Eintragungsdatum with new date value
$ws.Range("B9").Value = "01.01.2024"

# Update selection to match final state
$ws.Range("B10").Select()
